$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.827676240208877, 0.3952569169960474, 5.839057899901865, 1.648351648351648, 0.6575342465753425, 0.06863417982155114, 0.1531393568147014, 0.1055408970976253),
    @(1.253263707571802, 0.7315700619020821, 6.026457618814306, 2.188449848024316, 0.1098297638660077, 0.7529089664613279, 0.05107252298263534, 0.2109704641350211),
    @(4.647519582245431, 5.420666290231508, 6.032368808239333, 2.134146341463415, 0.4939626783754116, 0.3431708991077557, 0.0511770726714432, 4.931794333683106),
    @(14.83028720626632, 30.63877897117015, 5.977462028417443, 1.769371568029286, 0.4393190554640308, 2.200825309491059, 0, 0.0526592943654555),
    @(0.8337675872850442, 0.4512126339537508, 5.740922473012758, 1.828153564899452, 0.3841931942919868, 1.987662782727896, 0.2549719530851606, 5.628616517622304),
    @(0.5232862375719518, 0.2261164499717354, 6.136475208640157, 1.830384380719951, 0.1097694840834248, 2.537722908093278, 5.413687436159346, 8.175105485232068),
    @(8.407310704960835, 19.37677053824363, 6.237721021611002, 1.830384380719951, 0.1646542261251372, 2.747252747252747, 5.316973415132924, 8.258811152025251),
    @(0.6266318537859008, 0.1693958215697346, 6.185567010309279, 1.825928180158247, 0.7138934651290499, 0.6172839506172839, 0.3058103975535168, 0.2110817941952506),
    @(0.05227391531625718, 18.32579185520362, 5.937193326790972, 1.94765672550213, 0.6040637012630422, 0.2741603838245374, 5.660377358490567, 7.61904761904762),
    @(0.5227391531625719, 0.226628895184136, 5.795677799607073, 1.887941534713764, 0.2747252747252747, 1.517241379310345, 5.750636132315521, 7.871104067617539),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = 2 + $j  # Column B is index 2
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
